# Update the yearly income-statement database: drop the oldest period
# (1396/12), shift all remaining periods one column to the left, and
# append the newly published period (1401/12) together with its refreshed
# "read_price" figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "12 ماهه منتهی به ..." period headers -------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: "تاریخ انتشار" publish-date headers ----------------------------
$ws.Range("D9").Value = "1399-02-09 (10)"
$ws.Range("E9").Value = "1400-02-01 (11)"
$ws.Range("F9").Value = "1401-02-10 (12)"
$ws.Range("G9").Value = "1402-02-09 (10)"
$ws.Range("H9").Value = "1402-02-09 (2)"

# --- Row 11: فروش (Sales) ---------------------------------------------------
$ws.Range("D11").Value = 9937
$ws.Range("E11").Value = 12174
$ws.Range("F11").Value = 9859
$ws.Range("G11").Value = 11591
$ws.Range("H11").Value = 11925

# --- Row 12: بهای تمام شده کالای فروش رفته (COGS) ---------------------------
$ws.Range("D12").Value = -7161
$ws.Range("E12").Value = -7150
$ws.Range("F12").Value = -5383
$ws.Range("G12").Value = -6713
$ws.Range("H12").Value = -8108

# --- Row 13: سود (زیان) ناخالص (Gross profit) -------------------------------
$ws.Range("D13").Value = 2776
$ws.Range("E13").Value = 5024
$ws.Range("F13").Value = 4477
$ws.Range("G13").Value = 4877
$ws.Range("H13").Value = 3817

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی ------------------------------
$ws.Range("D14").Value = -718
$ws.Range("E14").Value = -672
$ws.Range("F14").Value = -472
$ws.Range("G14").Value = -593
$ws.Range("H14").Value = -747

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (هزینه استثنایی) -------------------
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = -63
$ws.Range("F15").Value = 11
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی -------------------------
$ws.Range("D16").Value = 20
$ws.Range("E16").Value = 21
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = "-"
$ws.Range("H16").Value = -10

# --- Row 17: سود (زیان) عملیاتی ---------------------------------------------
$ws.Range("D17").Value = 2078
$ws.Range("E17").Value = 4311
$ws.Range("F17").Value = 4015
$ws.Range("G17").Value = 4285
$ws.Range("H17").Value = 3060

# --- Row 18: هزینه های مالی --------------------------------------------------
$ws.Range("D18").Value = -198
$ws.Range("E18").Value = -94
$ws.Range("F18").Value = -13
$ws.Range("G18").Value = "-"
$ws.Range("H18").Value = -17

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ------------------------
$ws.Range("D19").Value = -5
$ws.Range("E19").Value = 599
$ws.Range("F19").Value = 376
$ws.Range("G19").Value = -34
$ws.Range("H19").Value = 43

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---------------
$ws.Range("D20").Value = 1874
$ws.Range("E20").Value = 4816
$ws.Range("F20").Value = 4379
$ws.Range("G20").Value = 4251
$ws.Range("H20").Value = 3086

# --- Row 21: مالیات ----------------------------------------------------------
$ws.Range("D21").Value = -460
$ws.Range("E21").Value = -601
$ws.Range("F21").Value = -559
$ws.Range("G21").Value = -690
$ws.Range("H21").Value = -151

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ------------------------------
$ws.Range("D22").Value = 1414
$ws.Range("E22").Value = 4215
$ws.Range("F22").Value = 3820
$ws.Range("G22").Value = 3560
$ws.Range("H22").Value = 2935

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (all "-", unchanged)
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"

# --- Row 24: سود (زیان) خالص -------------------------------------------------
$ws.Range("D24").Value = 1414
$ws.Range("E24").Value = 4215
$ws.Range("F24").Value = 3820
$ws.Range("G24").Value = 3560
$ws.Range("H24").Value = 2935

# --- Row 25: سود هر سهم پس از کسر مالیات (all 0, unchanged) ------------------
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# --- Row 26: سرمایه (Capital) — refreshed with updated read_price figures ---
$ws.Range("D26").Value = 4449
$ws.Range("E26").Value = 3508
$ws.Range("F26").Value = 3556
$ws.Range("G26").Value = 7580
$ws.Range("H26").Value = 5667

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه (all 0, unchanged) -------------
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
